$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.130.34"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.650.99"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").Value = "'216.53"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "'0.5227"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("E7").Value = "  -0.15%  "
$ws.Range("D8").Value = "'0.2613"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "'0.06334"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("E10").Value = "  -1.52%  "
$ws.Range("D11").Value = "'0.07694"
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "1.644.26"
$ws.Range("E12").Value = "  -1.26%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  -1.38%  "
$ws.Range("D15").Value = "'0.5596"
$ws.Range("E15").Value = "  +1.81%  "
$ws.Range("D16").Value = "0.0₅8241"
$ws.Range("E16").Value = "  +3.50%  "
$ws.Range("D17").Value = "'65.38"
$ws.Range("E17").Value = "  -1.53%  "
$ws.Range("D18").Value = "26.133.01"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("D20").Value = "'4.761"
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("D21").Value = "'189.15"
$ws.Range("E21").Value = "  +1.47%  "
$ws.Range("D22").Value = "'10.27"
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("D23").Value = "'6.232"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "'1.002"
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").Value = "'146.46"
$ws.Range("E25").Value = "  -1.90%  "
$ws.Range("D26").Value = "'7.470"
$ws.Range("E26").Value = "  -0.44%  "
$ws.Range("E27").Value = "  -2.80%  "
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "'1.394"
$ws.Range("E29").Value = "  +3.30%  "
$ws.Range("D30").Value = "'0.05919"
$ws.Range("E30").Value = "  -7.10%  "
$ws.Range("D31").Value = "'1.270"
$ws.Range("E31").Value = "  -0.39%  "
$ws.Range("D32").Value = "'3.441"
$ws.Range("E32").Value = "  -2.23%  "
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  +0.68%  "
$ws.Range("D35").Value = "'0.9900"
$ws.Range("E35").Value = "  -1.54%  "
$ws.Range("B36").Value = "MXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D36").Value = "'2.760"
$ws.Range("E36").Value = "  +0.75%  "
$ws.Range("B37").Value = "HuobiToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D37").Value = "'2.392"
$ws.Range("E37").Value = "  -0.70%  "
$ws.Range("D38").Value = "'0.5683"
$ws.Range("E38").Value = "  -5.46%  "
$ws.Range("D39").Value = "'0.01621"
$ws.Range("E39").Value = "  +0.32%  "
$ws.Range("D40").Value = "'0.8591"
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "'5.797"
$ws.Range("E41").Value = "  -5.75%  "
$ws.Range("E42").Value = "  -0.28%  "
$ws.Range("D43").Value = "1.030.72"
$ws.Range("E43").Value = "  -7.49%  "
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").Value = "1.799.57"
$ws.Range("E46").Value = "  -0.82%  "
$ws.Range("D47").Value = "'56.13"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.005"
$ws.Range("E48").Value = "  +0.11%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'8.095"
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("D50").Value = "'0.05191"
$ws.Range("D51").Value = "'0.4219"
$ws.Range("E51").Value = "  -0.43%  "
